# Applies row-data changes to "Hortaliza, Femacal de La Calera - Pepino dulce"
# sheet: rows 2-17 are reordered/updated by date-grouped record blocks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44424
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 75
$ws.Cells.Item(2, 11).Value = 18000
$ws.Cells.Item(2, 12).Value = 18000
$ws.Cells.Item(2, 13).Value = 18000
$ws.Cells.Item(2, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(2, 16).Value = 1200
$ws.Cells.Item(2, 17).Value = 15

$ws.Cells.Item(3, 4).Value = 44424
$ws.Cells.Item(3, 9).Value = 'Segunda'
$ws.Cells.Item(3, 10).Value = 50
$ws.Cells.Item(3, 11).Value = 12000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 12000
$ws.Cells.Item(3, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(3, 16).Value = 800
$ws.Cells.Item(3, 17).Value = 15

$ws.Cells.Item(4, 4).Value = 44235
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 14000
$ws.Cells.Item(4, 12).Value = 14000
$ws.Cells.Item(4, 13).Value = 14000
$ws.Cells.Item(4, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(4, 16).Value = 778
$ws.Cells.Item(4, 17).Value = 18

$ws.Cells.Item(5, 4).Value = 44235
$ws.Cells.Item(5, 9).Value = 'Segunda'
$ws.Cells.Item(5, 10).Value = 70
$ws.Cells.Item(5, 11).Value = 12000
$ws.Cells.Item(5, 12).Value = 12000
$ws.Cells.Item(5, 13).Value = 12000
$ws.Cells.Item(5, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(5, 16).Value = 667
$ws.Cells.Item(5, 17).Value = 18

$ws.Cells.Item(6, 4).Value = 44235
$ws.Cells.Item(6, 9).Value = 'Tercera'
$ws.Cells.Item(6, 10).Value = 60
$ws.Cells.Item(6, 11).Value = 10000
$ws.Cells.Item(6, 12).Value = 10000
$ws.Cells.Item(6, 13).Value = 10000
$ws.Cells.Item(6, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(6, 16).Value = 556
$ws.Cells.Item(6, 17).Value = 18

$ws.Cells.Item(7, 4).Value = 44991
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 75
$ws.Cells.Item(7, 11).Value = 13000
$ws.Cells.Item(7, 12).Value = 13000
$ws.Cells.Item(7, 13).Value = 13000
$ws.Cells.Item(7, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(7, 16).Value = 722
$ws.Cells.Item(7, 17).Value = 18

$ws.Cells.Item(8, 4).Value = 44991
$ws.Cells.Item(8, 9).Value = 'Segunda'
$ws.Cells.Item(8, 10).Value = 56
$ws.Cells.Item(8, 11).Value = 9000
$ws.Cells.Item(8, 12).Value = 9000
$ws.Cells.Item(8, 13).Value = 9000
$ws.Cells.Item(8, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(8, 16).Value = 500
$ws.Cells.Item(8, 17).Value = 18

$ws.Cells.Item(9, 4).Value = 44992
$ws.Cells.Item(9, 9).Value = 'Primera'
$ws.Cells.Item(9, 10).Value = 56
$ws.Cells.Item(9, 11).Value = 13000
$ws.Cells.Item(9, 12).Value = 13000
$ws.Cells.Item(9, 13).Value = 13000
$ws.Cells.Item(9, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(9, 16).Value = 722
$ws.Cells.Item(9, 17).Value = 18

$ws.Cells.Item(10, 4).Value = 44756
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 65
$ws.Cells.Item(10, 11).Value = 14000
$ws.Cells.Item(10, 12).Value = 14000
$ws.Cells.Item(10, 13).Value = 14000
$ws.Cells.Item(10, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(10, 16).Value = 933
$ws.Cells.Item(10, 17).Value = 15

$ws.Cells.Item(11, 4).Value = 44756
$ws.Cells.Item(11, 9).Value = 'Segunda'
$ws.Cells.Item(11, 10).Value = 68
$ws.Cells.Item(11, 11).Value = 12000
$ws.Cells.Item(11, 12).Value = 12000
$ws.Cells.Item(11, 13).Value = 12000
$ws.Cells.Item(11, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(11, 16).Value = 800
$ws.Cells.Item(11, 17).Value = 15

$ws.Cells.Item(12, 4).Value = 44238
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 90
$ws.Cells.Item(12, 11).Value = 13000
$ws.Cells.Item(12, 12).Value = 13000
$ws.Cells.Item(12, 13).Value = 13000
$ws.Cells.Item(12, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(12, 16).Value = 722
$ws.Cells.Item(12, 17).Value = 18

$ws.Cells.Item(13, 4).Value = 44238
$ws.Cells.Item(13, 9).Value = 'Segunda'
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 11000
$ws.Cells.Item(13, 12).Value = 11000
$ws.Cells.Item(13, 13).Value = 11000
$ws.Cells.Item(13, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(13, 16).Value = 611
$ws.Cells.Item(13, 17).Value = 18

$ws.Cells.Item(14, 4).Value = 44536
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 87
$ws.Cells.Item(14, 11).Value = 22000
$ws.Cells.Item(14, 12).Value = 22000
$ws.Cells.Item(14, 13).Value = 22000
$ws.Cells.Item(14, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(14, 16).Value = 1222
$ws.Cells.Item(14, 17).Value = 18

$ws.Cells.Item(15, 4).Value = 44536
$ws.Cells.Item(15, 9).Value = 'Segunda'
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 20000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 13).Value = 20000
$ws.Cells.Item(15, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(15, 16).Value = 1111
$ws.Cells.Item(15, 17).Value = 18

$ws.Cells.Item(16, 4).Value = 44242
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 13000
$ws.Cells.Item(16, 13).Value = 13000
$ws.Cells.Item(16, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(16, 16).Value = 722
$ws.Cells.Item(16, 17).Value = 18

$ws.Cells.Item(17, 4).Value = 44242
$ws.Cells.Item(17, 9).Value = 'Segunda'
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 10000
$ws.Cells.Item(17, 12).Value = 10000
$ws.Cells.Item(17, 13).Value = 10000
$ws.Cells.Item(17, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(17, 16).Value = 556
$ws.Cells.Item(17, 17).Value = 18
